$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "58.503.38"
$ws.Range("E2").Value = "  +0.31%  "
$ws.Range("D3").Value = "2.523.44"
$ws.Range("E3").Value = "  +1.58%  "
$ws.Range("D4").Value = "'0.997"
$ws.Range("E4").Value = "  -0.32%  "
$ws.Range("D5").Value = "'521.84"
$ws.Range("E5").Value = "  +0.22%  "
$ws.Range("D6").Value = "'133.23"
$ws.Range("E6").Value = "  -1.33%  "
$ws.Range("E7").Value = "  +0.35%  "
$ws.Range("D8").Value = "'0.563"
$ws.Range("E8").Value = "  +0.57%  "
$ws.Range("D9").Value = "2.522.75"
$ws.Range("E9").Value = "  +0.89%  "
$ws.Range("D10").Value = "'0.0977"
$ws.Range("E10").Value = "  -1.50%  "
$ws.Range("E11").Value = "  -1.61%  "
$ws.Range("E12").Value = "  -3.41%  "
$ws.Range("D13").Value = "'0.332"
$ws.Range("E13").Value = "  -2.48%  "
$ws.Range("D14").Value = "2.943.99"
$ws.Range("E14").Value = "  +0.52%  "
$ws.Range("D15").Value = "58.364.55"
$ws.Range("E15").Value = "  +0.20%  "
$ws.Range("D16").Value = "'22.14"
$ws.Range("E16").Value = "  -0.40%  "
$ws.Range("D17").Value = "'0.0000135"
$ws.Range("E17").Value = "  -0.66%  "
$ws.Range("D18").Value = "2.534.90"
$ws.Range("E18").Value = "  +1.49%  "
$ws.Range("D19").Value = "'10.67"
$ws.Range("E19").Value = "  -0.29%  "
$ws.Range("D20").Value = "'322.28"
$ws.Range("E20").Value = "  +0.30%  "
$ws.Range("D21").Value = "'4.16"
$ws.Range("E21").Value = "  -0.95%  "
$ws.Range("D22").Value = "'6.16"
$ws.Range("E22").Value = "  +6.83%  "
$ws.Range("D23").Value = "'0.997"
$ws.Range("E23").Value = "  -0.23%  "
$ws.Range("D24").Value = "'64.68"
$ws.Range("E24").Value = "  +0.33%  "
$ws.Range("E25").Value = "  -1.36%  "
$ws.Range("D26").Value = "'0.999"
$ws.Range("E26").Value = "  +0.23%  "
$ws.Range("D27").Value = "'0.160"
$ws.Range("E27").Value = "  -1.20%  "
$ws.Range("D28").Value = "'7.40"
$ws.Range("E28").Value = "  -0.15%  "
$ws.Range("D29").Value = "0.0₃0752"
$ws.Range("E29").Value = "  +0.09%  "
$ws.Range("D30").Value = "'168.58"
$ws.Range("E30").Value = "  -0.75%  "
$ws.Range("E31").Value = "  +1.08%  "
$ws.Range("B32").Value = "Aptos"
$ws.Range("C32").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D32").Value = "'6.30"
$ws.Range("E32").Value = "  -0.51%  "
$ws.Range("B33").Value = "Fetch.AI"
$ws.Range("C33").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D33").Value = "'1.19"
$ws.Range("E33").Value = "  -0.25%  "
$ws.Range("D34").Value = "'0.998"
$ws.Range("E34").Value = "  -0.06%  "
$ws.Range("D35").Value = "'0.999"
$ws.Range("E35").Value = "  +0.25%  "
$ws.Range("D36").Value = "'18.18"
$ws.Range("E36").Value = "  +0.10%  "
$ws.Range("E37").Value = "  -6.27%  "
$ws.Range("D38").Value = "'3.92"
$ws.Range("E38").Value = "  -3.05%  "
$ws.Range("E39").Value = "  +0.82%  "
$ws.Range("D40").Value = "'36.45"
$ws.Range("E40").Value = "  -0.64%  "
$ws.Range("D41").Value = "'0.770"
$ws.Range("E41").Value = "  -3.84%  "
$ws.Range("D42").Value = "'276.63"
$ws.Range("E42").Value = "  -0.88%  "
$ws.Range("D43").Value = "'3.45"
$ws.Range("E43").Value = "  -0.48%  "
$ws.Range("D44").Value = "'129.94"
$ws.Range("E44").Value = "  +4.59%  "
$ws.Range("E45").Value = "  -4.10%  "
$ws.Range("E46").Value = "  -0.43%  "
$ws.Range("D47").Value = "'0.0917"
$ws.Range("E47").Value = "  +0.70%  "
$ws.Range("D48").Value = "'0.0500"
$ws.Range("E48").Value = "  +1.62%  "
$ws.Range("D49").Value = "'17.70"
$ws.Range("E49").Value = "  -0.73%  "
$ws.Range("E50").Value = "  -0.08%  "
$ws.Range("D51").Value = "'16.91"
$ws.Range("E51").Value = "  -1.24%  "
